$d = $word.ActiveDocument

# The signature block previously read "       Yanhao Chen" under the
# "By: ___" line; the sender/signer is now Matt.
$d.Content.Find.Execute("Yanhao Chen", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Matt", 2)

# The signer's title changes from "President" to "Legal Attorney"
# (keep the "Its:  " label intact).
$d.Content.Find.Execute("President", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Legal Attorney", 2)
